# Avances Etiquetado Roboflow 5/29/2025
# Fill in the new tracking row (row 26) with the 29/5/2025 update and
# move the on-screen selection to where the user ended up after typing it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D26").Value = "29/5/2025"
$ws.Range("E26").Value = 90
$ws.Range("F26").Value = 476
$ws.Range("G26").Value = 0
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 650
$ws.Range("J26").Value = "N/A"

# Reflect the saved view state: scrolled so column E is leftmost, row 3 on
# top, with the active cell resting on I31 (just below the new row).
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 5
$ws.Range("I31").Select()
